# Weekly update: insert a new daily price record at the top of the
# Haba (Feria Lagunitas de Puerto Montt) data set, pushing the
# existing historical rows down by one and extending the table with
# the row that used to be last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above row 5 (first data row after the header
# and the 3 existing entries). Excel shifts rows 5:73 down to 6:74,
# copies formatting from the row above, and grows the used range /
# dimension to A1:R74 automatically.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44545
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 12000
$ws.Range("N5").Value = "`$/saco 25 kilos"
$ws.Range("O5").Value = "Región de La Araucanía"
$ws.Range("P5").Value = 480
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"
